# Generate Report for Handback
#
# Populates the "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns (F, G, H) on the zh-cn and de-de
# worksheets now that translations have come back, and flips the
# Status column from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# Cornflower blue (matches the workbook's existing HyperLink font color FF6495ED)
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) using Excel's BGR-packed long

function Set-HandbackRow {
    param(
        $ws,
        [int]$row,
        [string]$targetFileUrl,
        [string]$targetFileName,
        [string]$handbackFileUrl,
        [string]$handbackFileName,
        [string]$handbackDateTime
    )

    # F = Latest Target File (mirrors the "Source File Name" hyperlink in column A)
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $targetFileName
    $ws.Hyperlinks.Add($fCell, $targetFileUrl, "", "", $targetFileName) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # G = Latest Handback File (mirrors the "Latest Handoff File" hyperlink in column D)
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $handbackFileName
    $ws.Hyperlinks.Add($gCell, $handbackFileUrl, "", "", $handbackFileName) | Out-Null
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkColor

    # H = Latest Handback DateTime
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $handbackDateTime
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

Set-HandbackRow $wsZh 2 `
    "https://github.com/OpenLocalizationTest/oltest/blob/886c261894e544e64b451dfcc7343ee8902bd3d5/e2e/68452ea2-ae53-4303-89cb-ebba6981771f.md" `
    "68452ea2-ae53-4303-89cb-ebba6981771f.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b81bc256d4ec7c22efcea35a0ca655a2cfa348f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/68452ea2-ae53-4303-89cb-ebba6981771f.3be1799586f40b204af6cde91994eb9c0d84460e.zh-cn.xlf" `
    "68452ea2-ae53-4303-89cb-ebba6981771f.3be1799586f40b204af6cde91994eb9c0d84460e.zh-cn.xlf" `
    "2016-03-12 16:11:25"

Set-HandbackRow $wsZh 3 `
    "https://github.com/OpenLocalizationTest/oltest/blob/886c261894e544e64b451dfcc7343ee8902bd3d5/e2e/c8f6e92b-9339-465f-8035-71958be13724.md" `
    "c8f6e92b-9339-465f-8035-71958be13724.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b81bc256d4ec7c22efcea35a0ca655a2cfa348f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c8f6e92b-9339-465f-8035-71958be13724.6fb101080f3d1250941b2bd22220a16b0d2e6d7e.zh-cn.xlf" `
    "c8f6e92b-9339-465f-8035-71958be13724.6fb101080f3d1250941b2bd22220a16b0d2e6d7e.zh-cn.xlf" `
    "2016-03-12 16:11:25"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

Set-HandbackRow $wsDe 2 `
    "https://github.com/OpenLocalizationTest/oltest/blob/886c261894e544e64b451dfcc7343ee8902bd3d5/e2e/68452ea2-ae53-4303-89cb-ebba6981771f.md" `
    "68452ea2-ae53-4303-89cb-ebba6981771f.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14547cb3b0ceb8567fe53c07e156f4db9f22c5b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/68452ea2-ae53-4303-89cb-ebba6981771f.3be1799586f40b204af6cde91994eb9c0d84460e.de-de.xlf" `
    "68452ea2-ae53-4303-89cb-ebba6981771f.3be1799586f40b204af6cde91994eb9c0d84460e.de-de.xlf" `
    "2016-03-12 16:11:31"

Set-HandbackRow $wsDe 3 `
    "https://github.com/OpenLocalizationTest/oltest/blob/886c261894e544e64b451dfcc7343ee8902bd3d5/e2e/c8f6e92b-9339-465f-8035-71958be13724.md" `
    "c8f6e92b-9339-465f-8035-71958be13724.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14547cb3b0ceb8567fe53c07e156f4db9f22c5b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c8f6e92b-9339-465f-8035-71958be13724.6fb101080f3d1250941b2bd22220a16b0d2e6d7e.de-de.xlf" `
    "c8f6e92b-9339-465f-8035-71958be13724.6fb101080f3d1250941b2bd22220a16b0d2e6d7e.de-de.xlf" `
    "2016-03-12 16:11:31"
